$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as plain text to avoid Excel auto-converting
# numeric-looking strings (e.g. "26.20") into actual numbers, which
# would drop trailing zeros / introduce floating point artifacts.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '60.878.02'
$ws.Range("E2").Value = '  +3.72%  '
Set-TextValue $ws.Range("D3") '3.233.82'
$ws.Range("E3").Value = '  +2.19%  '
Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.13%  '
Set-TextValue $ws.Range("D5") '542.44'
$ws.Range("E5").Value = '  +2.45%  '
Set-TextValue $ws.Range("D6") '147.43'
$ws.Range("E6").Value = '  +5.42%  '
Set-TextValue $ws.Range("D7") '0.998'
$ws.Range("E7").Value = '  -0.14%  '
Set-TextValue $ws.Range("D8") '0.532'
$ws.Range("E8").Value = '  -1.06%  '
Set-TextValue $ws.Range("D9") '7.37'
$ws.Range("E9").Value = '  +1.11%  '
Set-TextValue $ws.Range("D10") '0.114'
$ws.Range("E10").Value = '  +2.59%  '
Set-TextValue $ws.Range("D11") '0.437'
$ws.Range("E11").Value = '  -0.23%  '
Set-TextValue $ws.Range("D12") '3.767.37'
$ws.Range("E12").Value = '  +1.56%  '
Set-TextValue $ws.Range("D14") '26.20'
$ws.Range("E14").Value = '  +1.61%  '
Set-TextValue $ws.Range("D15") '0.0000174'
$ws.Range("E15").Value = '  +2.34%  '
Set-TextValue $ws.Range("D16") '60.744.78'
$ws.Range("E16").Value = '  +3.42%  '
Set-TextValue $ws.Range("D17") '3.234.19'
$ws.Range("E17").Value = '  +1.79%  '
Set-TextValue $ws.Range("D18") '6.33'
$ws.Range("E18").Value = '  +1.28%  '
Set-TextValue $ws.Range("D19") '13.37'
$ws.Range("E19").Value = '  +3.01%  '
Set-TextValue $ws.Range("D20") '8.35'
$ws.Range("E20").Value = '  +3.14%  '
Set-TextValue $ws.Range("D21") '377.71'
$ws.Range("E21").Value = '  +0.40%  '
Set-TextValue $ws.Range("D23") '0.528'
$ws.Range("E23").Value = '  -0.14%  '
Set-TextValue $ws.Range("D24") '70.04'
$ws.Range("E24").Value = '  +0.48%  '
Set-TextValue $ws.Range("D25") '0.172'
$ws.Range("E25").Value = '  +2.56%  '
Set-TextValue $ws.Range("D26") '8.72'
$ws.Range("E26").Value = '  +4.87%  '
Set-TextValue $ws.Range("D27") '0.999'
$ws.Range("E27").Value = '  -0.09%  '
Set-TextValue $ws.Range("D28") '0.0₃0899'
$ws.Range("E28").Value = '  +4.22%  '
Set-TextValue $ws.Range("D29") '22.64'
$ws.Range("E29").Value = '  +1.08%  '
Set-TextValue $ws.Range("D30") '1.91'
$ws.Range("E30").Value = '  +1.42%  '
Set-TextValue $ws.Range("D31") '6.19'
$ws.Range("E31").Value = '  +2.64%  '
Set-TextValue $ws.Range("D32") '5.34'
$ws.Range("E32").Value = '  +3.62%  '
Set-TextValue $ws.Range("D33") '1.21'
$ws.Range("E33").Value = '  +5.51%  '
Set-TextValue $ws.Range("D34") '6.64'
$ws.Range("E34").Value = '  +5.13%  '
Set-TextValue $ws.Range("D35") '158.58'
$ws.Range("E35").Value = '  +1.09%  '
Set-TextValue $ws.Range("D36") '1.40'
$ws.Range("E36").Value = '  +4.39%  '
Set-TextValue $ws.Range("D37") '26.60'
$ws.Range("E37").Value = '  +6.24%  '
Set-TextValue $ws.Range("D38") '2.807.95'
$ws.Range("E38").Value = '  +4.38%  '
Set-TextValue $ws.Range("D41") '1.71'
$ws.Range("E41").Value = '  +1.28%  '
Set-TextValue $ws.Range("D42") '4.27'
$ws.Range("E42").Value = '  -0.38%  '
Set-TextValue $ws.Range("D44") '0.724'
$ws.Range("E44").Value = '  -0.02%  '
Set-TextValue $ws.Range("D46") '3.262.93'
$ws.Range("E46").Value = '  +1.74%  '
Set-TextValue $ws.Range("D47") '0.997'
$ws.Range("E47").Value = '  +1.93%  '
Set-TextValue $ws.Range("D48") '6.20'
$ws.Range("E48").Value = '  -0.20%  '
Set-TextValue $ws.Range("D49") '20.91'
$ws.Range("E49").Value = '  +4.22%  '
Set-TextValue $ws.Range("D50") '0.806'
$ws.Range("E50").Value = '  +7.62%  '
Set-TextValue $ws.Range("D51") '0.998'
$ws.Range("E51").Value = '  -0.14%  '

$ws.Range("E13").Value = '  -1.85%  '
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("E45").Value = '  +2.65%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D39") '0.0715'
$ws.Range("E39").Value = '  +3.09%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D40") '0.0315'
$ws.Range("E40").Value = '  +8.51%  '

